$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 16 de Abril de 2020 a las 08:52"

# Row 44
$ws.Range("B44").Value = 4161
$ws.Range("C44").Value = 397
$ws.Range("D44").Value = 186
$ws.Range("E44").Value = 3859
$ws.Range("F44").Value = 45
$ws.Range("G44").Value = 8
$ws.Range("H44").Value = 116

# Row 70
$ws.Range("B70").Value = 1341
$ws.Range("C70").Value = 46
$ws.Range("D70").Value = 240
$ws.Range("E70").Value = 1085
$ws.Range("F70").Value = 22
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 16

# Row 74
$ws.Range("A74").Value = "Lituania"
$ws.Range("B74").Value = 1128
$ws.Range("C74").Value = 37
$ws.Range("D74").Value = 178
$ws.Range("E74").Value = 920
$ws.Range("F74").Value = 14
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 30

# Row 75
$ws.Range("A75").Value = "Armenia"
$ws.Range("B75").Value = 1111
$ws.Range("C75").Value = 0
$ws.Range("D75").Value = 297
$ws.Range("E75").Value = 797
$ws.Range("F75").Value = 30
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 17

# Row 76
$ws.Range("A76").Value = "Bosnia y Herzegovina"
$ws.Range("B76").Value = 1110
$ws.Range("C76").Value = 0
$ws.Range("D76").Value = 253
$ws.Range("E76").Value = 816
$ws.Range("F76").Value = 4
$ws.Range("G76").Value = 0
$ws.Range("H76").Value = 41

# Row 77
$ws.Range("A77").Value = "Oman"
$ws.Range("B77").Value = 1019
$ws.Range("C77").Value = 109
$ws.Range("D77").Value = 176
$ws.Range("E77").Value = 839
$ws.Range("F77").Value = 3
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 4

# Row 78
$ws.Range("A78").Value = "Hong Kong"
$ws.Range("B78").Value = 1017
$ws.Range("C78").Value = 0
$ws.Range("D78").Value = 459
$ws.Range("E78").Value = 554
$ws.Range("F78").Value = 10
$ws.Range("G78").Value = 0
$ws.Range("H78").Value = 4

# Row 79
$ws.Range("A79").Value = "Republica de Macedonia"
$ws.Range("B79").Value = 974
$ws.Range("C79").Value = 0
$ws.Range("D79").Value = 98
$ws.Range("E79").Value = 831
$ws.Range("F79").Value = 15
$ws.Range("G79").Value = 0
$ws.Range("H79").Value = 45

# Row 105
$ws.Range("B105").Value = 399
$ws.Range("C105").Value = 0
$ws.Range("D105").Value = 82
$ws.Range("E105").Value = 314
$ws.Range("F105").Value = 4
$ws.Range("G105").Value = 0
$ws.Range("H105").Value = 3

# Row 110
$ws.Range("A110").Value = "Georgia"
$ws.Range("B110").Value = 336
$ws.Range("C110").Value = 30
$ws.Range("D110").Value = 74
$ws.Range("E110").Value = 259
$ws.Range("F110").Value = 6
$ws.Range("G110").Value = 0
$ws.Range("H110").Value = 3

# Row 111
$ws.Range("A111").Value = "Mauricio"
$ws.Range("B111").Value = 324
$ws.Range("C111").Value = 0
$ws.Range("D111").Value = 65
$ws.Range("E111").Value = 250
$ws.Range("F111").Value = 3
$ws.Range("G111").Value = 0
$ws.Range("H111").Value = 9

# Row 112
$ws.Range("A112").Value = "Senegal"
$ws.Range("B112").Value = 314
$ws.Range("C112").Value = 0
$ws.Range("D112").Value = 190
$ws.Range("E112").Value = 122
$ws.Range("F112").Value = 1
$ws.Range("G112").Value = 0
$ws.Range("H112").Value = 2

# Row 113
$ws.Range("B113").Value = 290
$ws.Range("C113").Value = 2
$ws.Range("D113").Value = 55
$ws.Range("E113").Value = 231
$ws.Range("F113").Value = 7
$ws.Range("G113").Value = 0
$ws.Range("H113").Value = 4

# Row 117
$ws.Range("B117").Value = 238
$ws.Range("C117").Value = 0
$ws.Range("D117").Value = 65
$ws.Range("E117").Value = 166
$ws.Range("F117").Value = 1
$ws.Range("G117").Value = 0
$ws.Range("H117").Value = 7

# Row 122
$ws.Range("B122").Value = 184
$ws.Range("C122").Value = 0
$ws.Range("D122").Value = 169
$ws.Range("E122").Value = 15
$ws.Range("F122").Value = 0
$ws.Range("G122").Value = 0
$ws.Range("H122").Value = 0

# Row 184
$ws.Range("A184").Value = "Suazilandia"
$ws.Range("B184").Value = 16
$ws.Range("C184").Value = 1
$ws.Range("D184").Value = 8
$ws.Range("E184").Value = 8
$ws.Range("F184").Value = 0
$ws.Range("G184").Value = 0
$ws.Range("H184").Value = 0

# Row 185
$ws.Range("A185").Value = "Dominica"
$ws.Range("B185").Value = 16
$ws.Range("C185").Value = 0
$ws.Range("D185").Value = 8
$ws.Range("E185").Value = 8
$ws.Range("F185").Value = 0
$ws.Range("G185").Value = 0
$ws.Range("H185").Value = 0
